$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra week columns (I, J) and the extra task rows (6, 7)
$ws.Columns("I:J").Delete() | Out-Null
$ws.Rows("6:7").Delete() | Out-Null

# Collapse the two-month header into a single month spanning all 3 week columns
$ws.Range("G2").UnMerge() | Out-Null
$ws.Range("G2").ClearContents() | Out-Null
$ws.Range("F2").Value = "Month 1"
$ws.Range("F2:H2").Merge() | Out-Null
$ws.Range("G2:H2").ClearFormats() | Out-Null

# Generic week labels
$ws.Range("F3").Value = "Week 1"
$ws.Range("G3").Value = "Week 2"
$ws.Range("H3").Value = "Week 3"

# Clear the now-unused start/end date helper cells
$ws.Range("D4").ClearContents() | Out-Null
$ws.Range("E4").ClearContents() | Out-Null
$ws.Range("D5").ClearContents() | Out-Null
$ws.Range("E5").ClearContents() | Out-Null
